$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing "Late" / "Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()

# Update the selected cell shown when the sheet is active.
$ws.Activate()
$ws.Range("I21").Select()
